$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the existing E/G values on rows 3 and 4 (F stays as-is).
#    Row3: E3 1E-8 -> 1E-7 ; G3 1E-4 -> 5.5E-4
#    Row4: E4 1E-8 -> 1E-6 ; G4 1E-4 -> 1E-3
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = 0.0000001
$ws.Range("G3").Value = 0.00055

$ws.Range("E4").Value = 0.000001
$ws.Range("G4").Value = 0.001

# ---------------------------------------------------------------------------
# 2. Fill in the new data (C, E, F, G, H, I) for rows 5-10. These rows
#    already had A/B/D populated; this adds the rest of the experiment grid,
#    repeating the same (E,F,G) triple pattern used in rows 2-4, with C
#    stepping 5 (rows 5-7) then 10 (rows 8-10).
# ---------------------------------------------------------------------------
$eVals = @(0.00000001, 0.0000001, 0.000001)
$fVals = @(0.000001, 0.00001, 0.0001)
$gVals = @(0.0001, 0.00055, 0.001)
$cVals = @(5, 5, 5, 10, 10, 10)

for ($i = 0; $i -lt 6; $i++) {
    $r = 5 + $i
    $pattern = $i % 3

    $ws.Range("C$r").Value = $cVals[$i]
    $ws.Range("E$r").Value = $eVals[$pattern]
    $ws.Range("E$r").NumberFormat = $ws.Range("E2").NumberFormat
    $ws.Range("F$r").Value = $fVals[$pattern]
    $ws.Range("F$r").NumberFormat = $ws.Range("F2").NumberFormat
    $ws.Range("G$r").Value = $gVals[$pattern]

    $ws.Range("H$r").Formula = "=B$r * C$r"
    $ws.Range("I$r").Formula = "=H$r/(SQRT(2.65 * 9.81 * (D$r^5)))"
}

# ---------------------------------------------------------------------------
# 3. Drop the now-unused rows 11-20 (data trimmed down to just rows 2-10).
# ---------------------------------------------------------------------------
$ws.Range("A11:N20").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 4. View-state tweaks: zoom to 140%, and move the live selection to J9.
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("J9").Select() | Out-Null
$excel.ActiveWindow.Zoom = 140

$wb.Application.Calculate() | Out-Null
